$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (leave A2, B2, F2 unchanged; update C2, D2, E2, G2)
$ws.Range("C2").Value2 = "C:\Users\InterviewRoom1\Phase1\Backend\emotions"
$ws.Range("D2").Value2 = "C:\Users\InterviewRoom1\Phase1\Videos\men.mp4"
$ws.Range("E2").Value2 = 22.0 / 86400.0
$ws.Range("G2").Value2 = "C:\Users\InterviewRoom1\Phase1\Backend\data-normalizer"

# Delete row 3 entirely (shifts nothing up since it's the last row)
$ws.Rows("3:3").Delete() | Out-Null

# Set selection to E3 as in the target
$ws.Range("E3").Select() | Out-Null

Write-Host "Done"
